$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.438.63"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "1.918.99"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.59"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4808"
$ws.Range("E7").Value = "  -0.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4053"
$ws.Range("E8").Value = "  -0.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08190"
$ws.Range("E9").Value = "  +1.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.36"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.927.65"
$ws.Range("E12").Value = "  +2.20%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.071"
$ws.Range("E13").Value = "  +1.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.236"
$ws.Range("E14").Value = "  +2.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.68"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06874"
$ws.Range("E16").Value = "  +2.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.59"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.009"
$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").Value = "29.445.57"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.675"
$ws.Range("E22").Value = "  +2.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.79"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.190"
$ws.Range("E24").Value = "  +1.70%  "

$ws.Range("D25").Value = "2.145.41"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.648"
$ws.Range("E26").Value = "  +8.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.63"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.06"
$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.107"
$ws.Range("E29").Value = "  +0.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.39"
$ws.Range("E30").Value = "  +1.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.014"
$ws.Range("E31").Value = "  -1.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09645"
$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.594"
$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.557"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.374"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06329"
$ws.Range("E36").Value = "  +3.94%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02284"
$ws.Range("E37").Value = "  +1.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.178"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.83"
$ws.Range("E39").Value = "  +5.93%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5944"
$ws.Range("E40").Value = "  +0.90%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.010"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1850"
$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.282"
$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.393"
$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.49"
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07529"
$ws.Range("E46").Value = "  -3.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5587"
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.952"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.66"
$ws.Range("E49").Value = "  +3.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.430"
$ws.Range("E50").Value = "  +3.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.07"
$ws.Range("E51").Value = "  -0.45%  "
